# Setting up demonstration for MEEM team
# Adds "biomass_observed" / "biomass_cutoff" columns to the "species" sheet,
# right after w_min and before ks, with per-species values (or "NA" for
# species that don't have an observed biomass / cutoff yet).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("species")
$ws.Activate()

# Insert two new blank columns at D:E (pushes ks..pisc from D:M to F:O).
$ws.Range("D1:E1").EntireColumn.Insert()

# Headers for the two new columns.
$ws.Range("D1").Value = "biomass_observed"
$ws.Range("E1").Value = "biomass_cutoff"

# predators
$ws.Range("D2").Value = 107.269812380122
$ws.Range("E2").Value = 25

# herbivores
$ws.Range("D3").Value = 33.793040863225002
$ws.Range("E3").Value = 25

# inverts (no observed biomass data yet)
$ws.Range("D4").Value = "NA"
$ws.Range("E4").Value = "NA"

# Resize all columns to fit their (now wider) contents.
$widths = @(9.6666666666667, 10.1666666666667, 11.1666666666667, 17, 13.8333333333333, 4.1666666666667, 2.6666666666667, 5.1666666666667, 1.3333333333333, 4.1666666666667, 5.3333333333333, 6.6666666666667, 10.8333333333333, 8.6666666666667, 5.3333333333333)
for ($i = 0; $i -lt $widths.Length; $i++) {
    $ws.Columns.Item($i + 1).ColumnWidth = $widths[$i]
}

# Restore the selected cell as it was left after the edit.
$ws.Range("H12").Select() | Out-Null
